$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14, shifting existing rows 14+ down by one.
$ws.Rows.Item(14).Insert()

# Populate the new row with the new entry (alphabetically between
# "E-filing exemption - Supreme Court" and "Emergency Order of Protection - Cook County").
$ws.Range("A14").Value = "Emergency Motion to Claim Exemption"
$ws.Range("B14").Value = "https://www.illinoislegalaid.org/legal-information/emergency-motion-claim-exemption"

# Match the visual style of the other URL cells in column B (Hyperlink cell style),
# without actually creating a clickable hyperlink relationship.
$ws.Range("B14").Style = "Hyperlink"

# Reflect the selection saved in the workbook after the edit.
$ws.Range("A20").Select()
